{"js": "// Office.js (Word JavaScript API) script.\n// Body of: async (context) => { ... }\n//\n// Implements the two visible content edits from the diff:\n//   1. Heading \"4. Models Used\" -> \"4. Methodology\"\n//   2. Trim the trailing sentence from the \"Visualization\" bullet so it\n//      reads \"...accumulation zones. \" instead of \"...accumulation zones.\n//      Volume overlays were added to check if the stock is under\n//      accumulation or distribution phases.\"\n\nconst body = context.document.body;\n\n// 1) \"4. Models Used\" -> \"4. Methodology\" (keep the \"4. \" prefix, swap the\n//    remaining wording so the run boundary stays close to the original).\nconst headingResults = body.search(\"Models Used\", { matchCase: true, matchWholeWord: false });\nheadingResults.load(\"items\");\nawait context.sync();\n\nif (headingResults.items.length > 0) {\n  headingResults.items[0].insertText(\"Methodology\", Word.InsertLocation.replace);\n}\n\n// 2) Remove the trailing \"Volume overlays...\" sentence from the\n//    Visualization bullet, leaving the preceding sentence (and its\n//    trailing space) untouched.\nconst volumeSentence =\n  \"Volume overlays were added to check if the stock is under accumulation or distribution phases.\";\nconst volumeResults = body.search(volumeSentence, { matchCase: true, matchWholeWord: false });\nvolumeResults.load(\"items\");\nawait context.sync();\n\nif (volumeResults.items.length > 0) {\n  volumeResults.items[0].insertText(\"\", Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) script.\n# $word / $d resolve to the running application / the open document.\n#\n# Implements the two visible content edits from the diff:\n#   1. Heading \"4. Models Used\" -> \"4. Methodology\"\n#   2. Trim the trailing sentence from the \"Visualization\" bullet so it\n#      reads \"...accumulation zones. \" instead of \"...accumulation zones.\n#      Volume overlays were added to check if the stock is under\n#      accumulation or distribution phases.\"\n\n$d = $word.ActiveDocument\n\n# 1) \"4. Models Used\" -> \"4. Methodology\" (keep the \"4. \" prefix, swap the\n#    remaining wording so the run boundary stays close to the original).\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"Models Used\"\n$find.Replacement.Text = \"Methodology\"\n$find.Execute(\n  $find.Text,\n  $false,\n  $false,\n  $false,\n  $false,\n  $false,\n  $true,\n  1,\n  $false,\n  $find.Replacement.Text,\n  2\n) | Out-Null\n\n# 2) Remove the trailing \"Volume overlays...\" sentence from the\n#    Visualization bullet, leaving the preceding sentence (and its\n#    trailing space) untouched.\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Replacement.ClearFormatting()\n$find2.Text = \"Volume overlays were added to check if the stock is under accumulation or distribution phases.\"\n$find2.Replacement.Text = \"\"\n$find2.Execute(\n  $find2.Text,\n  $false,\n  $false,\n  $false,\n  $false,\n  $false,\n  $true,\n  1,\n  $false,\n  $find2.Replacement.Text,\n  2\n) | Out-Null\n"}
